$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Tempo Decorrido" (elapsed time) values for rows 2-15 (column F),
# reflecting the newly added "expired register" filter timings for SMERP searches.
$values = @(
    "0 hours, 0 minutes, and 32.78 seconds",
    "0 hours, 0 minutes, and 13.49 seconds",
    "0 hours, 0 minutes, and 13.70 seconds",
    "0 hours, 0 minutes, and 29.90 seconds",
    "0 hours, 0 minutes, and 14.35 seconds",
    "0 hours, 0 minutes, and 30.49 seconds",
    "0 hours, 0 minutes, and 12.29 seconds",
    "0 hours, 0 minutes, and 14.68 seconds",
    "0 hours, 0 minutes, and 12.98 seconds",
    "0 hours, 0 minutes, and 13.38 seconds",
    "0 hours, 0 minutes, and 28.75 seconds",
    "0 hours, 0 minutes, and 13.66 seconds",
    "0 hours, 0 minutes, and 31.00 seconds",
    "0 hours, 0 minutes, and 34.66 seconds"
)

$row = 2
foreach ($value in $values) {
    $ws.Cells.Item($row, 6).Value = $value
    $row = $row + 1
}
